# Update the cryptos list with refreshed prices / volume percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as plain text. Some new values look like ordinary
# decimal numbers, so force those specific cells to Text format first to
# keep them stored as strings (matching the rest of the column) instead of
# letting Excel auto-convert them to numeric cells.
$textPriceCells = @("D5","D7","D10","D11","D20","D21","D22","D23","D25","D26","D29","D30","D31","D32","D35","D37","D39","D41","D42","D47")
foreach ($cellRef in $textPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "41.978.38"
$ws.Range("E2").Value = "  -0.55%  "

$ws.Range("D3").Value = "2.218.47"
$ws.Range("E3").Value = "  -1.41%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").Value = "242.63"
$ws.Range("E5").Value = "  -1.94%  "

$ws.Range("E6").Value = "  -0.30%  "

$ws.Range("D7").Value = "73.69"
$ws.Range("E7").Value = "  -1.02%  "

$ws.Range("E8").Value = "  +0.13%  "

$ws.Range("E9").Value = "  -0.75%  "

$ws.Range("D10").Value = "43.66"
$ws.Range("E10").Value = "  +6.21%  "

$ws.Range("D11").Value = "0.0958"
$ws.Range("E11").Value = "  +1.91%  "

$ws.Range("E12").Value = "  +0.23%  "

$ws.Range("E13").Value = "  +0.02%  "

$ws.Range("D14").Value = "2.552.07"
$ws.Range("E14").Value = "  -1.31%  "

$ws.Range("E15").Value = "  -1.83%  "

$ws.Range("E16").Value = "  -1.18%  "

$ws.Range("D17").Value = "2.235.02"
$ws.Range("E17").Value = "  -0.83%  "

$ws.Range("D18").Value = "41.901.56"
$ws.Range("E18").Value = "  -0.48%  "

$ws.Range("E19").Value = "  +12.55%  "

$ws.Range("D20").Value = "6.19"
$ws.Range("E20").Value = "  +1.00%  "

$ws.Range("D21").Value = "72.42"
$ws.Range("E21").Value = "  +0.75%  "

$ws.Range("D22").Value = "10.48"
$ws.Range("E22").Value = "  +33.28%  "

$ws.Range("D23").Value = "229.63"
$ws.Range("E23").Value = "  -0.93%  "

$ws.Range("E24").Value = "  -7.31%  "

# Row 25 and 26 swap: Cosmos <-> Dai (including link and data)
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.07%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "11.52"
$ws.Range("E26").Value = "  +3.35%  "

$ws.Range("E27").Value = "  +1.41%  "

$ws.Range("E28").Value = "  -1.36%  "

$ws.Range("D29").Value = "2.20"
$ws.Range("E29").Value = "  +2.19%  "

$ws.Range("D30").Value = "166.49"
$ws.Range("E30").Value = "  -1.81%  "

$ws.Range("D31").Value = "20.58"
$ws.Range("E31").Value = "  -0.44%  "

$ws.Range("D32").Value = "5.66"
$ws.Range("E32").Value = "  +15.97%  "

$ws.Range("E33").Value = "  -3.11%  "

$ws.Range("E34").Value = "  -0.23%  "

$ws.Range("D35").Value = "29.34"
$ws.Range("E35").Value = "  -2.40%  "

$ws.Range("E36").Value = "  -4.22%  "

$ws.Range("D37").Value = "4.29"
$ws.Range("E37").Value = "  -4.45%  "

$ws.Range("E38").Value = "  +0.40%  "

$ws.Range("D39").Value = "12.98"
$ws.Range("E39").Value = "  -4.36%  "

$ws.Range("E40").Value = "  -1.91%  "

$ws.Range("D41").Value = "65.46"
$ws.Range("E41").Value = "  +5.06%  "

$ws.Range("D42").Value = "5.66"
$ws.Range("E42").Value = "  -2.03%  "

$ws.Range("E43").Value = "  -1.65%  "

$ws.Range("E44").Value = "  +0.56%  "

$ws.Range("E45").Value = "  -3.36%  "

$ws.Range("E46").Value = "  +0.41%  "

$ws.Range("D47").Value = "2.42"
$ws.Range("E47").Value = "  +5.88%  "

$ws.Range("E48").Value = "  -0.37%  "

$ws.Range("E49").Value = "  -0.31%  "

$ws.Range("E50").Value = "  +0.81%  "

$ws.Range("D51").Value = "2.427.25"
$ws.Range("E51").Value = "  -1.28%  "
